$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged content)
$ws.Range("A1").Value = "Type"
$ws.Range("B1").Value = "Value"

# Existing categories reused as-is (indices 0,1,2 preserved)
$ws.Range("A2").Value = "Coal Briquettes"
$ws.Range("A3").Value = "Iron Ore"
$ws.Range("A4").Value = "Petroleum Gas"
$ws.Range("A5").Value = "Other Mineral"

# New categories - set in the exact order they first appear so the
# shared-string table is built in the same sequence as the target file.
$ws.Range("A10").Value = "Manganese Ore"
$ws.Range("A7").Value = "Copper Ore"
$ws.Range("A11").Value = "Aluminium Ore"
$ws.Range("A8").Value = "Zinc Ore"
$ws.Range("A9").Value = "Precious Metal Ore"
$ws.Range("A6").Value = "Crude Petroleum"

# Values (column B), ordered to match the sorted (descending) layout
$ws.Range("B2").Value = 109.095566316
$ws.Range("B3").Value = 87.851720001000004
$ws.Range("B4").Value = 67.328487191999997
$ws.Range("B5").Value = 8.2697793009999998
$ws.Range("B6").Value = 7.5440299729999998
$ws.Range("B7").Value = 4.8912432130000001
$ws.Range("B8").Value = 2.2128310440000001
$ws.Range("B9").Value = 1.562934727
$ws.Range("B10").Value = 1.388132554
$ws.Range("B11").Value = 0.86498837700000009

# Column A widens to fit the new, longer mineral names (column B is left
# untouched so its existing best-fit width/flag survive unchanged)
$ws.Columns.Item(1).ColumnWidth = 41

# Selection moves to reflect the author's last selected cell
$ws.Range("G7").Select()

# Refresh the sort state to cover the new data extent (A2:B11 / B1:B11)
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B1:B11"), 0, 2) | Out-Null
$ws.Sort.SetRange($ws.Range("A1:B11"))
$ws.Sort.Header = 1
$ws.Sort.Apply()
